# Generate Report for Handoff
#
# The handoff for this localization package finished, so the status flips
# from "In Translation" to "Ready for handoff"; the CI report is
# regenerated a minute later, so the "latest" timestamps advance too, and
# the (now longer) status column widens to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" --------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps advance by about a minute --------------------------------
# Overview's "Latest HO Xliff Generate Date" tracks the de-de handoff
# (the later of the two locales), matching de-de's "Latest Handoff Datetime".
$overview.Range("G2").Value = "2016-08-13 09:11:16"
$dede.Range("H2").Value = "2016-08-13 09:11:16"

# zh-cn's own "Latest Handoff Datetime" advances independently.
$zhcn.Range("H2").Value = "2016-08-13 09:11:09"

# --- Widen the Status columns so the longer text isn't truncated --------
# (numeric column indices -- letter-indexed Columns.Item("E") does not
# resolve correctly in this runtime)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
